# issue #5: add legislator_id, name, date into dataframe
#
# The stock-holdings sheet ("股票") gains three new trailing columns:
#   H = date, I = legislator_name, J = legislator_id
# Row 1 gets the new column headers, row 2 gets the values describing the
# single existing stock-holding record.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Headers (row 1) ----------------------------------------------------
$ws.Range("H1").Formula = '="date"'
$ws.Range("I1").Formula = '="legislator_name"'
$ws.Range("J1").Formula = '="legislator_id"'
$ws.Range("H1:J1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4163) | Out-Null

# Match the bold/bordered/centered look of the other header cells.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null

# --- Data (row 2) ---------------------------------------------------------
# "2012-04-27" looks like a date to Excel's normal cell-input parser, which
# would silently turn it into a date serial number instead of keeping it as
# text. Entering it as a formula result and pasting back only the computed
# value keeps it as literal text, same as the legislator name string.
$ws.Range("H2").Formula = '="2012-04-27"'
$ws.Range("I2").Formula = '="陳雪生"'
$ws.Range("H2:I2").Copy() | Out-Null
$ws.Range("H2:I2").PasteSpecial(-4163) | Out-Null

# legislator_id is plain numeric data.
$ws.Range("J2").Value = 1751

# Match the plain look of the rest of the data row.
$ws.Range("G2").Copy() | Out-Null
$ws.Range("H2:J2").PasteSpecial(-4122) | Out-Null
